# Apply the changes described by the commit diff:
#  - rename sheet "temperature_c" -> "temperature"
#  - move the active tab / selection from "genotype" to "temperature"
#  - update the lingering selection left behind on "genotype"

$wb = $excel.ActiveWorkbook

# 1) Update the selection that is left behind on the "genotype" sheet
#    (it was the active sheet before; its tabSelected flag goes away and
#    its cell selection changes from B15 to J19:J20 with active cell J19).
$wsGenotype = $wb.Worksheets.Item("genotype")
$wsGenotype.Activate() | Out-Null
$wsGenotype.Range("J19:J20").Select() | Out-Null

# 2) Rename "temperature_c" to "temperature"
$wsTemp = $wb.Worksheets.Item("temperature_c")
$wsTemp.Name = "temperature"

# 3) Make "temperature" the active/selected sheet (keeps its own Q34
#    selection, but now becomes tabSelected and the workbook's activeTab).
$wsTemp.Activate() | Out-Null
